$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting normalization -------------------------------------------
# The "Shop 'n Save" row's city cell (C2) was the one oddball cell using a
# Roboto font on a white fill. Normalize it to the plain Arial/theme text
# with no fill, then spread that look onto the two site rows ("Dylamato's
# Market" and "Main Street Farmers Market", currently rows 3 & 4) that
# will remain after the "Shop 'n Save" row below is deleted.
$ws.Range("C2").Interior.Pattern = -4142
$ws.Range("C2").Font.Name = "Arial"
$ws.Range("C2").Font.ThemeColor = 1

$ws.Range("C2").Copy()
$ws.Range("A3:E4").PasteSpecial(-4122)

# --- Data changes ---------------------------------------------------------
# Remove the "Shop 'n Save" / New Kensington row; everything below shifts up
$ws.Rows.Item(2).Delete()

# Insert a new row before the trailing note row, for the newly-added site
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).RowHeight = 15.75

$ws.Range("A4").Value = "East End Food Co-op"
$ws.Range("B4").Value = "7516 Meade St"
$ws.Range("C4").Value = "Pittsburgh"
$ws.Range("D4").Value = "PA"
$ws.Range("E4").Value = 15208

# Give the new row the plain (unthemed) Arial look, then highlight it
# yellow so the newly-added site stands out.
$ws.Range("A1").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$ws.Range("A4:E4").Interior.ColorIndex = 6
